$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" header column in H1, copying the formatting of the
# neighboring header cell (G1: bold, centered, bordered) so the new
# column matches the existing header row's look.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Populate the new column's data row with its value.
$ws.Range("H2").Value = 0
